$d = $word.ActiveDocument

# --- Paragraph 1: collapse the many split runs into one run with the same text ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.End = $r1.End - 1
$r1.Delete()
$r1.InsertAfter("This part of the document has stayed the same from version to version.  It shouldn't be shown if it doesn't change.  Otherwise, that would not be helping to compress the size of the changes.")

# --- Paragraph 3: collapse its runs into one run ---
$p2 = $d.Paragraphs.Item(3)
$r2 = $p2.Range
$r2.End = $r2.End - 1
$r2.Delete()
$r2.InsertAfter("This paragraph contains text that is outdated. It will be deleted in the near future.")

# --- Paragraph 5: collapse its runs into one run ---
$p3 = $d.Paragraphs.Item(5)
$r3 = $p3.Range
$r3.End = $r3.End - 1
$r3.Delete()
$r3.InsertAfter("It is important to spell check this document. On the other hand, a misspelled word isn't the end of the world. Nothing in the rest of this paragraph needs to be changed. Things can be added after it.")

# --- Remove the leftover hidden "_GoBack" bookmark from the last edit location ---
if ($d.Bookmarks.Item("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Append two new empty paragraphs at the very end of the document ---
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

$endRng2 = $d.Content
$endRng2.Collapse(0)
$endRng2.InsertParagraphAfter()
